$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1 : headers ----
$ws.Range("A1").Value = "Имя 1"

$ws.Range("B1").Value = "Партия 1"
$ws.Range("D1").Value = "Партия 2"
$ws.Range("F1").Value = "Партия 3"
$ws.Range("H1").Value = "Партия 4"
$ws.Range("J1").Value = "Партия 5"
$ws.Range("L1").Value = "Партия 6"
$ws.Range("N1").Value = "Партия 7"

$ws.Range("P1").Value = "Имя 2"
$ws.Range("Q1").Value = "Общий счет"

# Merge each "Партия n" header pair first ...
$ws.Range("B1:C1").Merge()
$ws.Range("D1:E1").Merge()
$ws.Range("F1:G1").Merge()
$ws.Range("H1:I1").Merge()
$ws.Range("J1:K1").Merge()
$ws.Range("L1:M1").Merge()
$ws.Range("N1:O1").Merge()

# ... then center them all, so every pair shares one style entry.
$ws.Range("B1:C1").HorizontalAlignment = -4108
$ws.Range("D1:E1").HorizontalAlignment = -4108
$ws.Range("F1:G1").HorizontalAlignment = -4108
$ws.Range("H1:I1").HorizontalAlignment = -4108
$ws.Range("J1:K1").HorizontalAlignment = -4108
$ws.Range("L1:M1").HorizontalAlignment = -4108
$ws.Range("N1:O1").HorizontalAlignment = -4108

# ---- Row 2 : data ----
$ws.Range("A2").Value = "Петя"
$ws.Range("B2").Value = 11
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 12
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 13
$ws.Range("G2").Value = 11
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = "Витя"

# Score column: typed as "3:1" which Excel first tries to read as a time
# (hence the h:mm number format) but keeps as literal text (quote-prefixed).
$ws.Range("Q2").Value = "'3:1"
$ws.Range("Q2").NumberFormat = "h:mm"

# ---- Column width ----
$ws.Columns("Q").ColumnWidth = 11.7109375

# ---- View / selection ----
$ws.Range("K17").Select()

$wb.Windows.Item(1).Left = -27630
$wb.Windows.Item(1).Top = 3525
